$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to text format while writing, so numeric-looking
# strings (e.g. "0.998") are not silently converted to numbers by Excel.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "51.513.13"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "2.983.01"
$ws.Range("E3").Value = "  +2.99%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "380.97"
$ws.Range("E5").Value = "  +4.08%  "

$ws.Range("D6").Value = "104.69"
$ws.Range("E6").Value = "  +2.93%  "

$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +3.15%  "

$ws.Range("D10").Value = "37.28"
$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  +2.43%  "

$ws.Range("D13").Value = "3.445.11"
$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("D14").Value = "18.44"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  +3.34%  "

$ws.Range("D16").Value = "2.983.84"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("D17").Value = "0.975"
$ws.Range("E17").Value = "  +6.85%  "

$ws.Range("D18").Value = "51.464.35"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").Value = "3.33"
$ws.Range("E19").Value = "  +4.36%  "

$ws.Range("D20").Value = "7.44"
$ws.Range("E20").Value = "  +4.62%  "

$ws.Range("D21").Value = "13.00"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("D23").Value = "69.33"
$ws.Range("E23").Value = "  +2.41%  "

$ws.Range("D24").Value = "262.83"
$ws.Range("E24").Value = "  +2.44%  "

$ws.Range("D25").Value = "2.89"

$ws.Range("D26").Value = "8.37"
$ws.Range("E26").Value = "  +21.16%  "

$ws.Range("D27").Value = "7.80"
$ws.Range("E27").Value = "  +27.96%  "

$ws.Range("E28").Value = "  +16.97%  "

$ws.Range("E29").Value = "  +2.71%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "26.01"
$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "9.92"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("D33").Value = "35.09"
$ws.Range("E33").Value = "  +3.97%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "2.08"
$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.04"
$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").Value = "0.0454"
$ws.Range("E36").Value = "  +8.91%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "3.05"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").Value = "17.16"
$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("E42").Value = "  +4.73%  "

$ws.Range("D43").Value = "125.48"
$ws.Range("E43").Value = "  +5.83%  "

$ws.Range("D44").Value = "21.82"
$ws.Range("E44").Value = "  +1.12%  "

$ws.Range("E45").Value = "  +21.30%  "

$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("D48").Value = "2.038.08"
$ws.Range("E48").Value = "  +1.68%  "

$ws.Range("D49").Value = "3.27"

$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +8.67%  "

$ws.Range("D51").Value = "1.30"
$ws.Range("E51").Value = "  +4.51%  "

# Restore default (General) formatting now that the text values are locked in.
$valueRange.ClearFormats()
